$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 298.66666
$ws.Range("I9").Value = 298
$ws.Range("K9").Value = 298
$ws.Range("M9").Value = -129
$ws.Range("H18").Value = 7403.2
$ws.Range("I18").Value = 9085.375
$ws.Range("K18").Value = 9085.375
$ws.Range("M18").Value = -8801.375
$ws.Range("H51").Value = 14853.889
$ws.Range("J51").Value = 9210.75
$ws.Range("L51").Value = 9210.75
$ws.Range("N51").Value = -10178.75
$ws.Range("H106").Value = 1000000
$ws.Range("I106").Value = 1000000
$ws.Range("K106").Value = 1000000
$ws.Range("M106").Value = -999369
$ws.Range("H112").Value = 5186.6665
$ws.Range("I112").Value = 399.25
$ws.Range("J112").Value = 5653.732
$ws.Range("K112").Value = 1197.75
$ws.Range("L112").Value = 16961.196
$ws.Range("M112").Value = -89.75
$ws.Range("N112").Value = -19177.196
$ws.Range("H113").Value = 80560760
$ws.Range("I113").Value = 123458540
$ws.Range("K113").Value = 123458540
$ws.Range("M113").Value = -123455286
$ws.Range("H137").Value = 2388.5715
$ws.Range("I137").Value = 2389.15
$ws.Range("K137").Value = 7167.450000000001
$ws.Range("M137").Value = -4617.450000000001
$ws.Range("H138").Value = 4231.5835
$ws.Range("J138").Value = 5622.933
$ws.Range("L138").Value = 16868.799
$ws.Range("N138").Value = -27148.799

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2810.6365
$ws.Range("I2").Value = 983.93335
$ws.Range("J2").Value = 6725
$ws.Range("K2").Value = 983.93335
$ws.Range("L2").Value = 6725
$ws.Range("M2").Value = -870.93335
$ws.Range("N2").Value = -6951
$ws.Range("H32").Value = 3504.22
$ws.Range("I32").Value = 2969.837
$ws.Range("J32").Value = 9649.625
$ws.Range("K32").Value = 2969.837
$ws.Range("L32").Value = 9649.625
$ws.Range("M32").Value = -2682.837
$ws.Range("N32").Value = -10223.625
$ws.Range("H74").Value = 38254.977
$ws.Range("I74").Value = 52891.594
$ws.Range("J74").Value = 4799.857
$ws.Range("K74").Value = 52891.594
$ws.Range("L74").Value = 4799.857
$ws.Range("M74").Value = -52017.594
$ws.Range("N74").Value = -6547.857
$ws.Range("H77").Value = 38254.977
$ws.Range("I77").Value = 52891.594
$ws.Range("J77").Value = 4799.857
$ws.Range("K77").Value = 264457.97
$ws.Range("L77").Value = 23999.285
$ws.Range("M77").Value = -260089.97
$ws.Range("N77").Value = -32735.285
$ws.Range("H97").Value = 16667012
$ws.Range("I97").Value = 431.5
$ws.Range("K97").Value = 431.5
$ws.Range("M97").Value = 64.5
$ws.Range("H116").Value = 2810.6365
$ws.Range("I116").Value = 983.93335
$ws.Range("J116").Value = 6725
$ws.Range("K116").Value = 983.93335
$ws.Range("L116").Value = 6725
$ws.Range("M116").Value = 1310.06665
$ws.Range("N116").Value = -11313
$ws.Range("H122").Value = 24259.3
$ws.Range("I122").Value = 30084.715
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 90254.145
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -87804.145
$ws.Range("N122").Value = -36900.001

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2810.6365
$ws.Range("I3").Value = 983.93335
$ws.Range("J3").Value = 6725
$ws.Range("K3").Value = 983.93335
$ws.Range("L3").Value = 6725
$ws.Range("M3").Value = -869.93335
$ws.Range("N3").Value = -6953
$ws.Range("H94").Value = 1653.6842
$ws.Range("I94").Value = 700.38464
$ws.Range("J94").Value = 3719.1667
$ws.Range("K94").Value = 700.38464
$ws.Range("L94").Value = 3719.1667
$ws.Range("M94").Value = -249.38464
$ws.Range("N94").Value = -4621.1667

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("H22").Value = 321.44446
$ws.Range("I22").Value = 324.125
$ws.Range("K22").Value = 324.125
$ws.Range("M22").Value = 25.875
$ws.Range("H60").Value = 16856.428
$ws.Range("J60").Value = 37497.5
$ws.Range("L60").Value = 37497.5
$ws.Range("N60").Value = -38519.5
$ws.Range("H105").Value = 5498521.5
$ws.Range("I105").Value = 7937598
$ws.Range("K105").Value = 7937598
$ws.Range("M105").Value = -7935851
$ws.Range("H122").Value = 2266.5881
$ws.Range("J122").Value = 3623
$ws.Range("L122").Value = 10869
$ws.Range("N122").Value = -15769
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2355112.2
$ws.Range("I5").Value = 3333684
$ws.Range("J5").Value = 6540
$ws.Range("K5").Value = 10001052
$ws.Range("L5").Value = 19620
$ws.Range("M5").Value = -10000940
$ws.Range("N5").Value = -19844
$ws.Range("H56").Value = 7076
$ws.Range("I56").Value = 7076
$ws.Range("K56").Value = 7076
$ws.Range("M56").Value = -6546
$ws.Range("H68").Value = 25004098
$ws.Range("J68").Value = 66674610
$ws.Range("L68").Value = 200023830
$ws.Range("N68").Value = -200025452
$ws.Range("H71").Value = 25004098
$ws.Range("J71").Value = 66674610
$ws.Range("L71").Value = 600071490
$ws.Range("N71").Value = -600079602
$ws.Range("H80").Value = 41671020
$ws.Range("I80").Value = 27781638
$ws.Range("J80").Value = 83339170
$ws.Range("K80").Value = 83344914
$ws.Range("L80").Value = 250017510
$ws.Range("M80").Value = -83343978
$ws.Range("N80").Value = -250019382
$ws.Range("H83").Value = 41671020
$ws.Range("I83").Value = 27781638
$ws.Range("J83").Value = 83339170
$ws.Range("K83").Value = 250034742
$ws.Range("L83").Value = 750052530
$ws.Range("M83").Value = -250030062
$ws.Range("N83").Value = -750061890
$ws.Range("H114").Value = 722.6667
$ws.Range("I114").Value = 593.5
$ws.Range("J114").Value = 787.25
$ws.Range("K114").Value = 1780.5
$ws.Range("L114").Value = 2361.75
$ws.Range("M114").Value = 1473.5
$ws.Range("N114").Value = -8869.75
$ws.Range("H135").Value = 2355112.2
$ws.Range("I135").Value = 3333684
$ws.Range("J135").Value = 6540
$ws.Range("K135").Value = 30003156
$ws.Range("L135").Value = 58860
$ws.Range("M135").Value = -30000621
$ws.Range("N135").Value = -63930

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 19900
$ws.Range("J54").Value = 19900
$ws.Range("L54").Value = 19900
$ws.Range("M54").Value = -20680
$ws.Range("H97").Value = 964.64703
$ws.Range("I97").Value = 908.9655
$ws.Range("K97").Value = 908.9655
$ws.Range("M97").Value = -412.9655

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5632.4
$ws.Range("I7").Value = 4581.885
$ws.Range("K7").Value = 4581.885
$ws.Range("M7").Value = -4469.885
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 1637.6666
$ws.Range("J22").Value = 2884.4285
$ws.Range("L22").Value = 2884.4285
$ws.Range("N22").Value = -3474.4285
$ws.Range("H27").Value = 1637.6666
$ws.Range("J27").Value = 2884.4285
$ws.Range("L27").Value = 2884.4285
$ws.Range("N27").Value = -3098.4285
$ws.Range("H126").Value = 5632.4
$ws.Range("I126").Value = 4581.885
$ws.Range("K126").Value = 13745.655
$ws.Range("M126").Value = -11275.655
$ws.Range("H132").Value = 8480407
$ws.Range("I132").Value = 16669596
$ws.Range("J132").Value = 8831.896000000001
$ws.Range("K132").Value = 50008788
$ws.Range("L132").Value = 26495.688
$ws.Range("M132").Value = -50006258
$ws.Range("N132").Value = -31555.688
$ws.Range("H136").Value = 7786.2705
$ws.Range("I136").Value = 1952.4736
$ws.Range("K136").Value = 5857.4208
$ws.Range("M136").Value = -3307.4208
$ws.Range("M16").ClearContents()

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8331
$ws.Range("I62").Value = 8994
$ws.Range("K62").Value = 8994
$ws.Range("M62").Value = -8370
$ws.Range("H65").Value = 8331
$ws.Range("I65").Value = 8994
$ws.Range("K65").Value = 44970
$ws.Range("M65").Value = -41850
$ws.Range("H122").Value = 98181.336
$ws.Range("I122").Value = 122966.18
$ws.Range("J122").Value = 7303.5557
$ws.Range("K122").Value = 368898.54
$ws.Range("L122").Value = 21910.6671
$ws.Range("M122").Value = -366448.54
$ws.Range("N122").Value = -26810.6671
$ws.Range("H132").Value = 9267427
$ws.Range("I132").Value = 13161402
$ws.Range("K132").Value = 39484206
$ws.Range("M132").Value = -39481676
$ws.Range("H136").Value = 16413451
$ws.Range("I136").Value = 24391336
$ws.Range("K136").Value = 73174008
$ws.Range("M136").Value = -73171458

